$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (closest achievable value given COM's pixel-snapped
# column-width rounding; target stored width is 15.42578125)
$ws.Columns.Item(1).ColumnWidth = 14.59

# Update cell values
$ws.Range("A1").Value = 0.22826733537822008
$ws.Range("A2").Value = -0.0059999999444642071
$ws.Range("A3").Value = -0.003999999952631228
$ws.Range("A4").Value = -0.0079999999123891996
$ws.Range("A5").Value = -0.0029999999523839804
$ws.Range("A6").Value = -0.0019999999503870214
$ws.Range("A7").Value = -0.0099999998796187484
$ws.Range("A8").Value = -0.0099999998799438217
$ws.Range("A9").Value = -0.0019999999534938695
$ws.Range("A10").Value = -0.0019999999559452419
$ws.Range("A11").Value = -0.0029999999471854721
$ws.Range("A12").Value = -0.0034999999431364337
$ws.Range("A13").Value = -0.003499999948290089
$ws.Range("A14").Value = -0.00799999991038014
$ws.Range("A15").Value = 0.04077141437882581
$ws.Range("A16").Value = -0.0019999999661890477
$ws.Range("A17").Value = -0.0019999999654203293
$ws.Range("A18").Value = -0.0039999999472879466
$ws.Range("A19").Value = -0.052472851261399001
$ws.Range("A20").Value = -0.0039999999585145218
$ws.Range("A21").Value = -0.0039999999581166179
$ws.Range("A22").Value = -0.0039999999577799983
$ws.Range("A23").Value = -0.0049999999408374407
$ws.Range("A24").Value = -0.019999999801833646
$ws.Range("A25").Value = -0.019999999799257928
$ws.Range("A26").Value = -0.0024999999438595211
$ws.Range("A27").Value = -0.0024999999409134332
$ws.Range("A28").Value = -0.0019999999324733508
$ws.Range("A29").Value = -0.0069999998792074081
$ws.Range("A30").Value = -0.059999999401492587
$ws.Range("A31").Value = 0.0079553214512628756
$ws.Range("A32").Value = -0.0099999998466344664
$ws.Range("A33").Value = -0.0039999999001345543
